$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update material rows with new plaster materials
$ws.Range("A2").Value = 717
$ws.Range("B2").Value = "Lime Cement Mortar (High Cement Ratio)"

$ws.Range("A3").Value = 718
$ws.Range("B3").Value = "Lime Cement Mortar (Low Cement Ratio)"

# Update the selected cell / active selection to B8
$ws.Range("B8").Select()
